# Generate Report for Handback
#
# - Marks both locales (zh-cn, de-de) as handed back / in sync with en-US.
# - Stamps the "Latest Handback DateTime" for each locale.
# - Adds the "Latest Target File" / "Latest Handback File" columns (F/G) for
#   each locale's two rows, with matching hyperlinks.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Status text: every cell that used to read "Ready for handoff" now
#     reads the handback status message. ---
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

# --- Latest Handback DateTime (column H): zh-cn synced a few seconds
#     before de-de. ---
$zhcn.Range("H2").Value = "2016-03-18 00:49:03"
$zhcn.Range("H3").Value = "2016-03-18 00:49:03"

$dede.Range("H2").Value = "2016-03-18 00:49:09"
$dede.Range("H3").Value = "2016-03-18 00:49:09"

function Add-HandbackColumns($ws) {
    # Grab the existing hyperlink targets for the source (.md) file and the
    # handoff (.xlf) file on each row so the new "Latest Target File" /
    # "Latest Handback File" links point at the same place.
    $aAddr = @{}
    $dAddr = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq "$A$2") { $aAddr[2] = $hl.Address }
        if ($addr -eq "$A$3") { $aAddr[3] = $hl.Address }
        if ($addr -eq "$D$2") { $dAddr[2] = $hl.Address }
        if ($addr -eq "$D$3") { $dAddr[3] = $hl.Address }
    }

    foreach ($row in 2, 3) {
        $aText = $ws.Range("A$row").Text
        $dText = $ws.Range("D$row").Text

        $fCell = $ws.Range("F$row")
        $fCell.Value = $aText
        $ws.Hyperlinks.Add($fCell, $aAddr[$row], [Type]::Missing, [Type]::Missing, $aText) | Out-Null
        $fCell.Font.Underline = 2
        $fCell.Font.Color = 0xED9564

        $gCell = $ws.Range("G$row")
        $gCell.Value = $dText
        $ws.Hyperlinks.Add($gCell, $dAddr[$row], [Type]::Missing, [Type]::Missing, $dText) | Out-Null
        $gCell.Font.Underline = 2
        $gCell.Font.Color = 0xED9564
    }
}

Add-HandbackColumns $zhcn
Add-HandbackColumns $dede
